# "Generate Report for Archive"
# - Status text "Ready for handoff" -> "In Translation" on every sheet that
#   shows it (Overview summary columns + each per-locale status column).
# - The narrower status text means the Status column(s) can shrink to fit,
#   so re-fit those columns' widths.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to [string] explicitly: Value2 can come back as a native
        # Boolean for "True"/"False" cells, and PowerShell's -eq coerces
        # the other operand using the *left-hand* operand's type, which
        # would otherwise make ($true -eq $oldStatus) spuriously true.
        $text = [string]$cell.Value2
        if ($text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# Re-fit the (now narrower) Status columns on every sheet.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)
